$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1054
$ws.Range("I6").Value = 283
$ws.Range("J6").Value = 1825
$ws.Range("K6").Value = 849
$ws.Range("L6").Value = 5475
$ws.Range("M6").Value = -737
$ws.Range("N6").Value = -5699

$ws.Range("H8").Value = 422.33334
$ws.Range("I8").Value = 106.8
$ws.Range("K8").Value = 320.4
$ws.Range("M8").Value = -181.4

$ws.Range("H11").Value = 198.3
$ws.Range("I11").Value = 198.3
$ws.Range("K11").Value = 198.3
$ws.Range("M11").Value = -58.30000000000001

$ws.Range("H17").Value = 993.6739
$ws.Range("J17").Value = 993.6739
$ws.Range("L17").Value = 2981.0217
$ws.Range("N17").Value = -3317.0217

$ws.Range("H101").Value = 1799.75
$ws.Range("I101").Value = 300
$ws.Range("J101").Value = 2299.6667
$ws.Range("K101").Value = 900
$ws.Range("L101").Value = 6899.000100000001
$ws.Range("M101").Value = 722
$ws.Range("N101").Value = -10143.0001

$ws.Range("H113").Value = 21667.727
$ws.Range("I113").Value = 23484.5
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 23484.5
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -20230.5
$ws.Range("N113").Value = -10008

$ws.Range("H137").Value = 1443.8846
$ws.Range("I137").Value = 1099.3
$ws.Range("J137").Value = 2592.5
$ws.Range("K137").Value = 3297.9
$ws.Range("L137").Value = 7777.5
$ws.Range("M137").Value = -747.8999999999996
$ws.Range("N137").Value = -12877.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7752604
$ws.Range("I2").Value = 7752604
$ws.Range("K2").Value = 7752604
$ws.Range("M2").Value = -7752491

$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H32").Value = 3708.1968
$ws.Range("I32").Value = 2126.0833
$ws.Range("K32").Value = 2126.0833
$ws.Range("M32").Value = -1839.0833

$ws.Range("H74").Value = 724.6
$ws.Range("I74").Value = 724.6
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 724.6
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 149.4
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 724.6
$ws.Range("I77").Value = 724.6
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 3623
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 745
$ws.Range("N77").ClearContents()

$ws.Range("H116").Value = 7752604
$ws.Range("I116").Value = 7752604
$ws.Range("K116").Value = 7752604
$ws.Range("M116").Value = -7750310

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7752604
$ws.Range("I3").Value = 7752604
$ws.Range("K3").Value = 7752604
$ws.Range("M3").Value = -7752490

$ws.Range("H86").Value = 134872.67
$ws.Range("I86").Value = 1262.7273
$ws.Range("J86").Value = 502300
$ws.Range("K86").Value = 1262.7273
$ws.Range("L86").Value = 502300
$ws.Range("M86").Value = -139.7273
$ws.Range("N86").Value = -504546

$ws.Range("H89").Value = 134872.67
$ws.Range("I89").Value = 1262.7273
$ws.Range("J89").Value = 502300
$ws.Range("K89").Value = 6313.636500000001
$ws.Range("L89").Value = 2511500
$ws.Range("M89").Value = -697.6365000000005
$ws.Range("N89").Value = -2522732

$ws.Range("H134").Value = 9041.737999999999
$ws.Range("I134").Value = 9230.794
$ws.Range("J134").Value = 8238.25
$ws.Range("K134").Value = 27692.382
$ws.Range("L134").Value = 24714.75
$ws.Range("M134").Value = -25157.382
$ws.Range("N134").Value = -29784.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 657.75
$ws.Range("I22").Value = 299.14285
$ws.Range("J22").Value = 1159.8
$ws.Range("K22").Value = 299.14285
$ws.Range("L22").Value = 1159.8
$ws.Range("M22").Value = 50.85714999999999
$ws.Range("N22").Value = -1859.8

$ws.Range("H94").Value = 984.2727
$ws.Range("I94").Value = 1003
$ws.Range("K94").Value = 1003
$ws.Range("M94").Value = -552

$ws.Range("H107").Value = 366.72
$ws.Range("I107").Value = 319.5238
$ws.Range("J107").Value = 614.5
$ws.Range("K107").Value = 319.5238
$ws.Range("L107").Value = 614.5
$ws.Range("M107").Value = 1600.4762
$ws.Range("N107").Value = -4454.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 629.6
$ws.Range("J26").Value = 282.66666
$ws.Range("L26").Value = 847.9999799999999
$ws.Range("N26").Value = -1423.99998

$ws.Range("H131").Value = 809.14435
$ws.Range("I131").Value = 318.14285
$ws.Range("J131").Value = 847.3333
$ws.Range("K131").Value = 954.4285500000001
$ws.Range("L131").Value = 2541.9999
$ws.Range("M131").Value = 4085.57145
$ws.Range("N131").Value = -12621.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2338.75
$ws.Range("I80").Value = 2200
$ws.Range("J80").Value = 2477.5
$ws.Range("K80").Value = 2200
$ws.Range("L80").Value = 2477.5
$ws.Range("M80").Value = -1202
$ws.Range("N80").Value = -4473.5

$ws.Range("H83").Value = 2338.75
$ws.Range("I83").Value = 2200
$ws.Range("J83").Value = 2477.5
$ws.Range("K83").Value = 11000
$ws.Range("L83").Value = 12387.5
$ws.Range("M83").Value = -6008
$ws.Range("N83").Value = -22371.5

$ws.Range("H113").Value = 1485.5714
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2727
$ws.Range("I16").Value = 3644.6
$ws.Range("J16").Value = 891.8
$ws.Range("K16").Value = 3644.6
$ws.Range("L16").Value = 891.8
$ws.Range("M16").Value = -3474.6
$ws.Range("N16").Value = -1231.8

$ws.Range("H55").Value = 478.7143
$ws.Range("I55").Value = 367
$ws.Range("J55").Value = 562.5
$ws.Range("K55").Value = 367
$ws.Range("L55").Value = 562.5
$ws.Range("M55").Value = -194
$ws.Range("N55").Value = -908.5

$ws.Range("H110").Value = 20000
$ws.Range("J110").Value = 20000
$ws.Range("L110").Value = 20000
$ws.Range("N110").Value = -28180

$ws.Range("H122").Value = 4236.9375
$ws.Range("I122").Value = 1980.5454
$ws.Range("K122").Value = 5941.6362
$ws.Range("M122").Value = -3491.6362

$ws.Range("H136").Value = 3552.92
$ws.Range("I136").Value = 1906.5834
$ws.Range("K136").Value = 5719.7502
$ws.Range("M136").Value = -3169.7502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 29703.334
$ws.Range("J70").Value = 29703.334
$ws.Range("L70").Value = 29703.334
$ws.Range("N70").Value = -30333.334

$ws.Range("H73").Value = 29703.334
$ws.Range("J73").Value = 29703.334
$ws.Range("L73").Value = 29703.334
$ws.Range("N73").Value = -31887.334

$ws.Range("H81").Value = 480
$ws.Range("I81").Value = 480
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 960
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 101
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 480
$ws.Range("I84").Value = 480
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 4800
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 504
$ws.Range("N84").ClearContents()

$ws.Range("H121").Value = 69999
$ws.Range("J121").Value = 69999
$ws.Range("L121").Value = 69999
$ws.Range("N121").Value = -73493

$ws.Range("H123").Value = 47550
$ws.Range("J123").Value = 47550
$ws.Range("L123").Value = 47550
$ws.Range("N123").Value = -57350

$ws.Range("H136").Value = 20578864
$ws.Range("I136").Value = 30866768
$ws.Range("J136").Value = 3055.2222
$ws.Range("K136").Value = 92600304
$ws.Range("L136").Value = 9165.6666
$ws.Range("M136").Value = -92597754
$ws.Range("N136").Value = -14265.6666

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 79093.63
$ws.Range("J141").Value = 79093.63
$ws.Range("L141").Value = 79093.63
